# Updated about section validation
#
# 1. "About" -> "value" on the Profile sheet changes from "OrangeHRM" to
#    "OrangeHRM OS 5.7".
# 2. The "url" hyperlink on the Login sheet (A2) is removed, and the cell's
#    formatting is restored to match the rest of the sheet (no more blue/
#    underlined hyperlink styling).

$wb = $excel.ActiveWorkbook
$login = $wb.Worksheets.Item("Login")
$profile = $wb.Worksheets.Item("Profile")

# Remove the hyperlink on Login!A2 and restore its (non-link) formatting by
# copying the plain format used by the rest of the sheet.
$linkCell = $login.Range("A2")
$linkCell.Hyperlinks.Delete()
$login.Range("B2").Copy()
$linkCell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Update the About/value text on the Profile sheet.
$profile.Range("B2").Value = "OrangeHRM OS 5.7"
